$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '47.325.42'
$ws.Range("E2").Value = '  +0.88%  '

$ws.Range("D3").Value = '2.490.72'
$ws.Range("E3").Value = '  +0.02%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.03'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.32%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.32'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.27%  '

$ws.Range("E7").Value = '  +0.43%  '

$ws.Range("E8").Value = '  -0.26%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.534'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.44%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.70'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.32%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0810'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.14%  '

$ws.Range("E12").Value = '  +0.24%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.37'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.61%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.11'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.99%  '

$ws.Range("D15").Value = '2.881.36'

$ws.Range("D16").Value = '2.490.61'
$ws.Range("E16").Value = '  -1.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.845'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.49%  '

$ws.Range("D18").Value = '47.230.61'
$ws.Range("E18").Value = '  +0.65%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.84'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.85%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.61'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.59%  '

$ws.Range("D21").Value = '0.0₃0933'
$ws.Range("E21").Value = '  -0.05%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.69'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +14.37%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.31'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.29%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '245.43'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.94%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.57'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.16%  '

$ws.Range("E26").Value = '  +0.20%  '

$ws.Range("E27").Value = '  -1.41%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.27'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.55%  '

$ws.Range("E29").Value = '  +0.10%  '

$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.60'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.88%  '

$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.136'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.27%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.66'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.58%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.81'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.37%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.34'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.89%  '

$ws.Range("E35").Value = '  +0.80%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.01'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.09%  '

$ws.Range("E37").Value = '  +1.76%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.66'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.82%  '

$ws.Range("E39").Value = '  -0.89%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '23.02'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.61%  '

$ws.Range("E41").Value = '  -0.25%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.23'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.04%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '117.96'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.12%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0296'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.16%  '

$ws.Range("D45").Value = '1.986.62'
$ws.Range("E45").Value = '  +1.14%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.03'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.78%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.01'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.82%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.10'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.30%  '

$ws.Range("E49").Value = '  -0.51%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.10'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.31%  '

$ws.Range("E51").Value = '  +4.63%  '
